# Increase font sizes throughout the resume document.
# Mapping (in points, i.e. Range.Font.Size):
#   16 -> 18   (name header)
#    9 -> 10   (contact info / body text / bullets / descriptions)
#   12 -> 13   (section headers)
#   10 -> 11   (overview paragraph)
#   11 -> 12   (job titles)
# Every other size is left untouched (defensive default: +1pt).

$d = $word.ActiveDocument
$count = $d.Paragraphs.Count

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $pStart = $p.Range.Start
    $pEnd = $p.Range.End

    # Exclude the trailing paragraph mark so we only touch the
    # run-level w:rPr/w:sz, not the paragraph mark's w:pPr/w:rPr/w:sz.
    $rangeEnd = $pEnd - 1
    if ($rangeEnd -le $pStart) {
        continue
    }

    $r = $d.Range($pStart, $rangeEnd)
    $sz = $r.Font.Size

    if ($sz -eq 16) {
        $r.Font.Size = 18
    } elseif ($sz -eq 9) {
        $r.Font.Size = 10
    } elseif ($sz -eq 12) {
        $r.Font.Size = 13
    } elseif ($sz -eq 10) {
        $r.Font.Size = 11
    } elseif ($sz -eq 11) {
        $r.Font.Size = 12
    }
}
